$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.585.42"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.645.33"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "3.117.99"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "65.458.87"
$ws.Range("D17").Value = "2.634.81"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "159.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0255"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0253"
$ws.Range("E49").Value = "  +10.36%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0990"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.93%  "
